$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "Quantidade Aprovada" column (N) next to the existing
# --- "Quantidade Autorizada" column (M) in the requisition table. ---

# Row 5 label for the new column (label block is merged N5:N9, mirrors M5:M9)
$ws.Range("N5").Value = "Quantidade Aprovada"
[void]$ws.Range("N5:N9").Merge()
# Materialize the (otherwise-empty) merged cells below the label with the
# same style as N5 so they round-trip as explicit cells like the rest of
# the merge blocks on this sheet.
$fillStyle = $ws.Range("N5").Style
$ws.Range("N6").Style = $fillStyle
$ws.Range("N7").Style = $fillStyle
$ws.Range("N8").Style = $fillStyle
$ws.Range("N9").Style = $fillStyle

# Row 10 template placeholders: M10 becomes the "authorized" placeholder,
# and the new N10 cell gets the "approved" placeholder that used to live
# in M10.
$ws.Range("M10").Value = "{product.quantityAuthorized}"
$ws.Range("N10").Value = "{product.quantityApproved}"

# Extend the row 1-4 label merges (right-hand info block) to cover the new
# column N as well.
[void]$ws.Range("M1:N1").Merge()
[void]$ws.Range("M2:N2").Merge()
[void]$ws.Range("M3:N3").Merge()
[void]$ws.Range("M4:N4").Merge()

# Extend the footer merge on row 13 to cover the new column N.
[void]$ws.Range("G13:M13").UnMerge()
[void]$ws.Range("G13:N13").Merge()

# Cursor position left where the author's last save recorded it.
$ws.Range("Q20").Select() | Out-Null
